# Kiosk workbook cleanup: drop the unused "Theater Bar" price sheet/table
# and rename the remaining sheet to the shorter "Angebot".

$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# Remove the "Theater Bar" sheet (and its Table2) entirely.
$wb.Worksheets.Item("Theater Bar").Delete() | Out-Null

# Rename the remaining sheet from "KinoKiosk Angebot 2023" to "Angebot".
$wb.Worksheets.Item("KinoKiosk Angebot 2023").Name = "Angebot"

# Make the remaining sheet active and leave the selection on F42.
$ws = $wb.Worksheets.Item("Angebot")
$ws.Activate() | Out-Null
$ws.Range("F42").Select() | Out-Null

$excel.DisplayAlerts = $true
